$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 12.421
$ws.Range("E7").Value = 12.28470000000001
$ws.Range("D8").Value = -8.711199999999986
$ws.Range("B12").Value = 5.632799999999995
$ws.Range("D12").Value = -7.956499999999997
$ws.Range("D14").Value = -8.508399999999996
$ws.Range("E19").Value = 13.0717
$ws.Range("E21").Value = 12.69890000000001
$ws.Range("D22").Value = -8.127999999999993
$ws.Range("E24").Value = 12.91629999999999
